# Edit "MES 01" worksheet: fill in zero values for previously-empty day
# total cells (A30, A32, ..., A64) and update the saved view/selection so
# that row 47 is at the top and A42 is the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")

# Rows whose column-A cell should now contain an explicit 0 value.
$rows = 30,32,34,36,38,40,42,44,46,48,50,52,54,56,58,60,62,64

foreach ($r in $rows) {
    $ws.Range("A$r").Value = 0
}

# Activate the sheet and set the view so that A47 is the top-left visible
# cell, and A42 is the current selection.
$ws.Activate()
$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
